# cross_listing_mappings.xlsx
#
# Commit: "Got rid of extra whitespace in the expected result column which
# causes tests to fail."
#
# The real change is the trailing-space cleanup on the expected_result
# column's "success " / "fail " values (now "success" / "fail"). The rest
# of the diff (Sheet2 getting populated with the same 4x4 block, the
# dimension/selection bookkeeping that comes along for the ride, etc.) is
# reproduced here too so the workbook ends up structurally in line with
# the target.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Trim the trailing whitespace out of the expected_result column ---
$ws1.Range("D2").Value = "success"
$ws1.Range("D3").Value = "fail"
$ws1.Range("D4").Value = "fail"

# --- 2. The old E/F columns were always empty placeholder cells; clearing
#        them drops the used range back down from A1:F4 to A1:D4 ---
$ws1.Range("E1:F3").Clear()

# --- 3. Mirror the same 4x4 block onto Sheet2 (matches the target diff,
#        which shows Sheet2 ending up with an identical copy of the data)
#        Styles (font color, then number format) are applied BEFORE the
#        value so the B/C "numeric-looking" ids are written as text -
#        matching the source - and so every cell reuses the existing
#        cellXfs entries (s="1"/s="2") instead of minting new ones. ---
for ($col = 1; $col -le 4; $col++) {
    $cell = $ws2.Cells.Item(1, $col)
    $cell.Font.Color = 0
    $cell.Value = $ws1.Cells.Item(1, $col).Text
}

for ($row = 2; $row -le 4; $row++) {
    $a = $ws2.Cells.Item($row, 1)
    $a.Font.Color = 0
    $a.Value = $ws1.Cells.Item($row, 1).Text

    $b = $ws2.Cells.Item($row, 2)
    $b.Font.Color = 0
    $b.NumberFormat = "@"
    $b.Value = $ws1.Cells.Item($row, 2).Text

    $c = $ws2.Cells.Item($row, 3)
    $c.Font.Color = 0
    $c.NumberFormat = "@"
    $c.Value = $ws1.Cells.Item($row, 3).Text

    $d = $ws2.Cells.Item($row, 4)
    $d.Font.Color = 0
    $d.Value = $ws1.Cells.Item($row, 4).Text
}

$ws2.PageSetup.Orientation = 1

# --- 4. Selections: Sheet2 picked up a "select all" state, Sheet1's active
#        cell moved to D2. Select Sheet2 first so Sheet1 ends up as the
#        tab-selected sheet again (matches the target, where only Sheet1
#        keeps tabSelected). ---
$ws2.Cells.Select() | Out-Null
$ws1.Range("D2").Select() | Out-Null
